# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-27
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 3
    8  = 1
    9  = 2
    10 = 3
    11 = 2
    12 = 0
    13 = 5
    14 = 3
    15 = 2
    16 = 1
    17 = 0
    18 = 4
    19 = 5
    20 = 1
    21 = 4
    22 = 5
    23 = 3
    24 = 1
    25 = 5
    26 = 3
    27 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
